$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(4).Delete()
$ws.Rows(4).AutoFit()
$ws.Rows(6).AutoFit()
$ws.Rows(11).AutoFit()
